$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Delete row 19 (the USA - USL CHAMPIONSHIP game), shifting row 20 up to row 19
$ws.Rows(19).Delete()

# Step 2: Update odds values in rows 2-18 (unaffected by the row deletion, since row 19 was below them)
# Row 2
$ws.Range("AC2").Value = 12
$ws.Range("AH2").Value = 9
$ws.Range("AM2").Value = 23
$ws.Range("BD2").Value = 151
$ws.Range("H2").Value = 3.5
$ws.Range("K2").Value = 2.25
$ws.Range("L2").Value = 2.75
$ws.Range("Q2").Value = 1.8
$ws.Range("R2").Value = 2
$ws.Range("Y2").Value = 11

# Row 3
$ws.Range("AA3").Value = 26
$ws.Range("AB3").Value = 34
$ws.Range("AC3").Value = 9
$ws.Range("AE3").Value = 13
$ws.Range("AG3").Value = 251
$ws.Range("AI3").Value = 11
$ws.Range("AJ3").Value = 9
$ws.Range("AP3").Value = 26
$ws.Range("AT3").Value = 2.75
$ws.Range("AX3").Value = 12
$ws.Range("G3").Value = 3.4
$ws.Range("I3").Value = 2.25
$ws.Range("K3").Value = 2.1
$ws.Range("L3").Value = 2.88
$ws.Range("M3").Value = 1.06
$ws.Range("N3").Value = 10
$ws.Range("O3").Value = 1.33
$ws.Range("P3").Value = 3.4
$ws.Range("Q3").Value = 2.07
$ws.Range("R3").Value = 1.83
$ws.Range("S3").Value = 1.4
$ws.Range("T3").Value = 2.75
$ws.Range("U3").Value = 1.8
$ws.Range("V3").Value = 1.95
$ws.Range("W3").Value = 10
$ws.Range("X3").Value = 17

# Row 4
$ws.Range("M4").Value = 1.05
$ws.Range("N4").Value = 11
$ws.Range("O4").Value = 1.25
$ws.Range("P4").Value = 4
$ws.Range("Q4").Value = 1.8
$ws.Range("R4").Value = 2

# Row 5
$ws.Range("AA5").Value = 15
$ws.Range("AD5").Value = 7
$ws.Range("AL5").Value = 26
$ws.Range("AQ5").Value = 34
$ws.Range("G5").Value = 1.91
$ws.Range("I5").Value = 3.75
$ws.Range("W5").Value = 9.5

# Row 8
$ws.Range("G8").Value = 2.2

# Row 10
$ws.Range("Q10").Value = 2.2
$ws.Range("R10").Value = 1.65

# Row 12
$ws.Range("AC12").Value = 11.5
$ws.Range("AD12").Value = 6.8
$ws.Range("AG12").Value = 350
$ws.Range("AH12").Value = 8
$ws.Range("AK12").Value = 17
$ws.Range("AL12").Value = 14.5
$ws.Range("AM12").Value = 23
$ws.Range("AO12").Value = 19.5
$ws.Range("AQ12").Value = 100
$ws.Range("AT12").Value = 2.87
$ws.Range("AV12").Value = 55
$ws.Range("AW12").Value = 3.85
$ws.Range("AX12").Value = 9.5
$ws.Range("AY12").Value = 17
$ws.Range("AZ12").Value = 35
$ws.Range("BA12").Value = 60
$ws.Range("BB12").Value = 200
$ws.Range("G12").Value = 3.7
$ws.Range("H12").Value = 3.45
$ws.Range("I12").Value = 1.9
$ws.Range("J12").Value = 4
$ws.Range("K12").Value = 2.15
$ws.Range("L12").Value = 2.47
$ws.Range("M12").Value = 1.01
$ws.Range("N12").Value = 11.5
$ws.Range("O12").Value = 1.23
$ws.Range("P12").Value = 3.35
$ws.Range("Q12").Value = 1.7
$ws.Range("R12").Value = 1.93
$ws.Range("S12").Value = 1.33
$ws.Range("T12").Value = 3.14
$ws.Range("U12").Value = 1.6
$ws.Range("V12").Value = 2.05
$ws.Range("W12").Value = 13
$ws.Range("X12").Value = 23

# Row 13
$ws.Range("AA13").Value = 25
$ws.Range("AC13").Value = 9.5
$ws.Range("AD13").Value = 8.25
$ws.Range("AE13").Value = 11.75
$ws.Range("AK13").Value = 18.5
$ws.Range("AN13").Value = 5.8
$ws.Range("AO13").Value = 17
$ws.Range("AT13").Value = 3.55
$ws.Range("AU13").Value = 6.4
$ws.Range("AW13").Value = 4.3
$ws.Range("AX13").Value = 9.25
$ws.Range("AY13").Value = 13.5
$ws.Range("BA13").Value = 40
$ws.Range("G13").Value = 3.35
$ws.Range("H13").Value = 3.95
$ws.Range("I13").Value = 1.9
$ws.Range("J13").Value = 3.65
$ws.Range("K13").Value = 2.42
$ws.Range("L13").Value = 2.35
$ws.Range("N13").Value = 9.5
$ws.Range("P13").Value = 4.85
$ws.Range("R13").Value = 2.55
$ws.Range("T13").Value = 3.55
$ws.Range("V13").Value = 2.57
$ws.Range("X13").Value = 23

# Row 14
$ws.Range("AC14").Value = 9
$ws.Range("AD14").Value = 8.75
$ws.Range("AF14").Value = 60
$ws.Range("AH14").Value = 20
$ws.Range("AL14").Value = 60
$ws.Range("AM14").Value = 50
$ws.Range("AR14").Value = 37
$ws.Range("AT14").Value = 3.25
$ws.Range("AW14").Value = 7.9
$ws.Range("AX14").Value = 35
$ws.Range("BA14").Value = 200
$ws.Range("G14").Value = 1.45
$ws.Range("H14").Value = 4.35
$ws.Range("I14").Value = 6.3
$ws.Range("J14").Value = 1.9
$ws.Range("K14").Value = 2.42
$ws.Range("L14").Value = 5.8
$ws.Range("N14").Value = 9
$ws.Range("P14").Value = 4.4
$ws.Range("Q14").Value = 1.53
$ws.Range("R14").Value = 2.35
$ws.Range("T14").Value = 3.25
$ws.Range("U14").Value = 1.7
$ws.Range("V14").Value = 2.05
$ws.Range("Y14").Value = 8
$ws.Range("Z14").Value = 10.75

# Row 15
$ws.Range("AA15").Value = 11.25
$ws.Range("AC15").Value = 27
$ws.Range("AD15").Value = 9.75
$ws.Range("AL15").Value = 32
$ws.Range("AN15").Value = 4.2
$ws.Range("AO15").Value = 7.7
$ws.Range("AP15").Value = 11.75
$ws.Range("AQ15").Value = 21
$ws.Range("AR15").Value = 32
$ws.Range("AS15").Value = 90
$ws.Range("AT15").Value = 3.9
$ws.Range("AW15").Value = 7.1
$ws.Range("G15").Value = 1.62
$ws.Range("I15").Value = 4.45
$ws.Range("J15").Value = 2.05
$ws.Range("K15").Value = 2.57
$ws.Range("L15").Value = 4.25
$ws.Range("P15").Value = 5.8
$ws.Range("Q15").Value = 1.34
$ws.Range("R15").Value = 3
$ws.Range("S15").Value = 1.21
$ws.Range("T15").Value = 3.9
$ws.Range("U15").Value = 1.39
$ws.Range("V15").Value = 2.77
$ws.Range("W15").Value = 13.5
$ws.Range("X15").Value = 12

# Row 16
$ws.Range("AA16").Value = 23
$ws.Range("AC16").Value = 10
$ws.Range("AD16").Value = 6
$ws.Range("AE16").Value = 13
$ws.Range("AG16").Value = 201
$ws.Range("AH16").Value = 8.5
$ws.Range("AM16").Value = 29
$ws.Range("AT16").Value = 2.75
$ws.Range("AY16").Value = 23
$ws.Range("BA16").Value = 67
$ws.Range("BB16").Value = 151
$ws.Range("H16").Value = 3.2
$ws.Range("I16").Value = 2.4
$ws.Range("J16").Value = 3.5
$ws.Range("L16").Value = 3.1
$ws.Range("M16").Value = 1.06
$ws.Range("N16").Value = 10
$ws.Range("O16").Value = 1.29
$ws.Range("P16").Value = 3.5
$ws.Range("Q16").Value = 1.98
$ws.Range("R16").Value = 1.88
$ws.Range("S16").Value = 1.4
$ws.Range("T16").Value = 2.75
$ws.Range("U16").Value = 1.7
$ws.Range("V16").Value = 2.05
$ws.Range("W16").Value = 10

# Row 17
$ws.Range("AD17").Value = 7.5
$ws.Range("U17").Value = 1.36
$ws.Range("V17").Value = 3
$ws.Range("Z17").Value = 23

# Row 18
$ws.Range("Q18").Value = 1.8
$ws.Range("R18").Value = 2

# Step 3: Update the new row 19 (previously row 20, the VENEZUELA - LIGA FUTVE game) with its tweaked values
$ws.Range("H19").Value = 3.35
$ws.Range("I19").Value = 5.3
$ws.Range("J19").Value = 2.22
$ws.Range("L19").Value = 5.4
$ws.Range("N19").Value = 7
$ws.Range("U19").Value = 2.02
$ws.Range("V19").Value = 1.62
$ws.Range("AA19").Value = 15.5
$ws.Range("AC19").Value = 7.8
$ws.Range("AD19").Value = 6.7
$ws.Range("AE19").Value = 19
$ws.Range("AF19").Value = 110
$ws.Range("AH19").Value = 12
$ws.Range("AJ19").Value = 17.5
$ws.Range("AL19").Value = 65
$ws.Range("AM19").Value = 70
$ws.Range("AN19").Value = 3.35
$ws.Range("AQ19").Value = 28
$ws.Range("AU19").Value = 7.9
$ws.Range("AV19").Value = 80
$ws.Range("AW19").Value = 6.6
$ws.Range("AY19").Value = 37
$ws.Range("AZ19").Value = 200
$ws.Range("BA19").Value = 250
